# Update cryptos list values (price & 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '70.997.99'
$cell.Style = 'Normal'
$cell = $ws.Range("E2")
$cell.NumberFormat = '@'
$cell.Value = '  +0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '3.849.91'
$cell.Style = 'Normal'
$cell = $ws.Range("E3")
$cell.NumberFormat = '@'
$cell.Value = '  +1.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D4")
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range("E4")
$cell.NumberFormat = '@'
$cell.Value = '  -0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '707.02'
$cell.Style = 'Normal'
$cell = $ws.Range("E5")
$cell.NumberFormat = '@'
$cell.Value = '  +1.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D6")
$cell.NumberFormat = '@'
$cell.Value = '172.80'
$cell.Style = 'Normal'
$cell = $ws.Range("E6")
$cell.NumberFormat = '@'
$cell.Value = '  +0.10%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D7")
$cell.NumberFormat = '@'
$cell.Value = '3.848.59'
$cell.Style = 'Normal'
$cell = $ws.Range("E7")
$cell.NumberFormat = '@'
$cell.Value = '  +1.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E8")
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D9")
$cell.NumberFormat = '@'
$cell.Value = '0.527'
$cell.Style = 'Normal'
$cell = $ws.Range("E9")
$cell.NumberFormat = '@'
$cell.Value = '  -0.48%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E10")
$cell.NumberFormat = '@'
$cell.Value = '  -0.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E11")
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E12")
$cell.NumberFormat = '@'
$cell.Value = '  -0.50%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '0.0000257'
$cell.Style = 'Normal'
$cell = $ws.Range("E13")
$cell.NumberFormat = '@'
$cell.Value = '  -0.69%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '36.73'
$cell.Style = 'Normal'
$cell = $ws.Range("E14")
$cell.NumberFormat = '@'
$cell.Value = '  +0.79%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D15")
$cell.NumberFormat = '@'
$cell.Value = '4.498.94'
$cell.Style = 'Normal'
$cell = $ws.Range("E15")
$cell.NumberFormat = '@'
$cell.Value = '  +1.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D16")
$cell.NumberFormat = '@'
$cell.Value = '3.834.66'
$cell.Style = 'Normal'
$cell = $ws.Range("E16")
$cell.NumberFormat = '@'
$cell.Value = '  +1.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '71.020.61'
$cell.Style = 'Normal'
$cell = $ws.Range("E17")
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D18")
$cell.NumberFormat = '@'
$cell.Value = '7.21'
$cell.Style = 'Normal'
$cell = $ws.Range("E18")
$cell.NumberFormat = '@'
$cell.Value = '  -0.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E19")
$cell.NumberFormat = '@'
$cell.Value = '  +0.72%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D20")
$cell.NumberFormat = '@'
$cell.Value = '17.37'
$cell.Style = 'Normal'
$cell = $ws.Range("E20")
$cell.NumberFormat = '@'
$cell.Value = '  -2.93%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '497.05'
$cell.Style = 'Normal'
$cell = $ws.Range("E21")
$cell.NumberFormat = '@'
$cell.Value = '  +2.91%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D22")
$cell.NumberFormat = '@'
$cell.Value = '10.65'
$cell.Style = 'Normal'
$cell = $ws.Range("E22")
$cell.NumberFormat = '@'
$cell.Value = '  -3.77%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D23")
$cell.NumberFormat = '@'
$cell.Value = '0.717'
$cell.Style = 'Normal'
$cell = $ws.Range("E23")
$cell.NumberFormat = '@'
$cell.Value = '  +0.44%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D24")
$cell.NumberFormat = '@'
$cell.Value = '85.57'
$cell.Style = 'Normal'
$cell = $ws.Range("E24")
$cell.NumberFormat = '@'
$cell.Value = '  +1.33%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E25")
$cell.NumberFormat = '@'
$cell.Value = '  +1.82%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D26")
$cell.NumberFormat = '@'
$cell.Value = '10.66'
$cell.Style = 'Normal'
$cell = $ws.Range("E26")
$cell.NumberFormat = '@'
$cell.Value = '  +1.66%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D27")
$cell.NumberFormat = '@'
$cell.Value = '12.20'
$cell.Style = 'Normal'
$cell = $ws.Range("E27")
$cell.NumberFormat = '@'
$cell.Value = '  -1.99%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E28")
$cell.NumberFormat = '@'
$cell.Value = '  -3.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E29")
$cell.NumberFormat = '@'
$cell.Value = '  +2.97%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D30")
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range("E30")
$cell.NumberFormat = '@'
$cell.Value = '  -0.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D31")
$cell.NumberFormat = '@'
$cell.Value = '7.53'
$cell.Style = 'Normal'
$cell = $ws.Range("E31")
$cell.NumberFormat = '@'
$cell.Value = '  -0.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D32")
$cell.NumberFormat = '@'
$cell.Value = '2.27'
$cell.Style = 'Normal'
$cell = $ws.Range("E32")
$cell.NumberFormat = '@'
$cell.Value = '  -0.92%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D33")
$cell.NumberFormat = '@'
$cell.Value = '29.47'
$cell.Style = 'Normal'
$cell = $ws.Range("E33")
$cell.NumberFormat = '@'
$cell.Value = '  -0.10%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E34")
$cell.NumberFormat = '@'
$cell.Value = '  -2.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D35")
$cell.NumberFormat = '@'
$cell.Value = '9.19'
$cell.Style = 'Normal'
$cell = $ws.Range("E35")
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D36")
$cell.NumberFormat = '@'
$cell.Value = '3.806.32'
$cell.Style = 'Normal'
$cell = $ws.Range("E36")
$cell.NumberFormat = '@'
$cell.Value = '  +1.43%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E37")
$cell.NumberFormat = '@'
$cell.Value = '  +0.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E38")
$cell.NumberFormat = '@'
$cell.Value = '  +0.24%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D39")
$cell.NumberFormat = '@'
$cell.Value = '2.37'
$cell.Style = 'Normal'
$cell = $ws.Range("E39")
$cell.NumberFormat = '@'
$cell.Value = '  +7.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '6.04'
$cell.Style = 'Normal'
$cell = $ws.Range("E40")
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E41")
$cell.NumberFormat = '@'
$cell.Value = '  +5.93%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E42")
$cell.NumberFormat = '@'
$cell.Value = '  -3.52%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E43")
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E44")
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E45")
$cell.NumberFormat = '@'
$cell.Value = '  -3.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '163.66'
$cell.Style = 'Normal'
$cell = $ws.Range("E46")
$cell.NumberFormat = '@'
$cell.Value = '  +0.80%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E47")
$cell.NumberFormat = '@'
$cell.Value = '  -0.37%  '
$cell.Style = 'Normal'

$cell = $ws.Range("E48")
$cell.NumberFormat = '@'
$cell.Value = '  +0.94%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D49")
$cell.NumberFormat = '@'
$cell.Value = '415.48'
$cell.Style = 'Normal'

$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '0.299'
$cell.Style = 'Normal'
$cell = $ws.Range("E50")
$cell.NumberFormat = '@'
$cell.Value = '  -1.42%  '
$cell.Style = 'Normal'

$cell = $ws.Range("D51")
$cell.NumberFormat = '@'
$cell.Value = '8.62'
$cell.Style = 'Normal'
$cell = $ws.Range("E51")
$cell.NumberFormat = '@'
$cell.Value = '  +0.29%  '
$cell.Style = 'Normal'

